# The "벨트" (Belt) product row for rep A / region 가 (row 5) was removed
# from the consolidated sales sheet. Deleting the whole row shifts every
# row below it up by one, shrinks the used range from A1:G11 to A1:G10,
# and drops the now-unused "벨트" shared string automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(5).Delete()
